# Auto-generated edit script applying numeric corrections to leve profit sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Cells.Item(17, 8).Value = 4250390.5
$ws.Cells.Item(17, 10).Value = 4420386
$ws.Cells.Item(17, 12).Value = 13261158
$ws.Cells.Item(17, 14).Value = -13261494
$ws.Cells.Item(62, 8).Value = 2279.3333
$ws.Cells.Item(62, 9).Value = 2559.1428
$ws.Cells.Item(62, 10).Value = 1300
$ws.Cells.Item(62, 11).Value = 2559.1428
$ws.Cells.Item(62, 12).Value = 1300
$ws.Cells.Item(62, 13).Value = -1935.1428
$ws.Cells.Item(62, 14).Value = -2548
$ws.Cells.Item(65, 8).Value = 2279.3333
$ws.Cells.Item(65, 9).Value = 2559.1428
$ws.Cells.Item(65, 10).Value = 1300
$ws.Cells.Item(65, 11).Value = 12795.714
$ws.Cells.Item(65, 12).Value = 6500
$ws.Cells.Item(65, 13).Value = -9675.714
$ws.Cells.Item(65, 14).Value = -12740
$ws.Cells.Item(132, 8).Value = 4632269.5
$ws.Cells.Item(132, 9).Value = 2609.7046
$ws.Cells.Item(132, 10).Value = 25002774
$ws.Cells.Item(132, 11).Value = 7829.1138
$ws.Cells.Item(132, 12).Value = 75008322
$ws.Cells.Item(132, 13).Value = -5299.1138
$ws.Cells.Item(132, 14).Value = -75013382
$ws.Cells.Item(141, 8).Value = 1363.5358
$ws.Cells.Item(141, 9).Value = 1128.2916
$ws.Cells.Item(141, 10).Value = 2775
$ws.Cells.Item(141, 11).Value = 3384.8748
$ws.Cells.Item(141, 12).Value = 8325
$ws.Cells.Item(141, 13).Value = 1795.1252
$ws.Cells.Item(141, 14).Value = -18685

$ws = $wb.Worksheets.Item("ARM")
$ws.Cells.Item(32, 8).Value = 4166.06
$ws.Cells.Item(32, 9).Value = 4021.5918
$ws.Cells.Item(32, 10).Value = 11245
$ws.Cells.Item(32, 11).Value = 4021.5918
$ws.Cells.Item(32, 12).Value = 11245
$ws.Cells.Item(32, 13).Value = -3734.5918
$ws.Cells.Item(32, 14).Value = -11819
$ws.Cells.Item(45, 8).Value = 1749934
$ws.Cells.Item(45, 9).Value = 2393997.8
$ws.Cells.Item(45, 11).Value = 2393997.8
$ws.Cells.Item(45, 13).Value = -2393620.8
$ws.Cells.Item(61, 8).Value = 2187.75
$ws.Cells.Item(61, 9).Value = 2193.1143
$ws.Cells.Item(61, 11).Value = 2193.1143
$ws.Cells.Item(61, 13).Value = -1981.1143
$ws.Cells.Item(122, 8).Value = 7074.3335
$ws.Cells.Item(122, 9).Value = 8754.666999999999
$ws.Cells.Item(122, 10).Value = 2033.3334
$ws.Cells.Item(122, 11).Value = 26264.001
$ws.Cells.Item(122, 12).Value = 6100.0002
$ws.Cells.Item(122, 13).Value = -23814.001
$ws.Cells.Item(122, 14).Value = -11000.0002
$ws.Cells.Item(136, 8).Value = 2187.75
$ws.Cells.Item(136, 9).Value = 2193.1143
$ws.Cells.Item(136, 11).Value = 6579.342900000001
$ws.Cells.Item(136, 13).Value = -4029.342900000001

$ws = $wb.Worksheets.Item("BSM")
$ws.Cells.Item(20, 8).Value = 2418.3572
$ws.Cells.Item(20, 9).Value = 2340.7273
$ws.Cells.Item(20, 10).Value = 2703
$ws.Cells.Item(20, 11).Value = 2340.7273
$ws.Cells.Item(20, 12).Value = 2703
$ws.Cells.Item(20, 13).Value = -2093.7273
$ws.Cells.Item(20, 14).Value = -3197
$ws.Cells.Item(107, 8).Value = 1239.3889
$ws.Cells.Item(107, 9).Value = 1109.1333
$ws.Cells.Item(107, 10).Value = 1890.6666
$ws.Cells.Item(107, 11).Value = 1109.1333
$ws.Cells.Item(107, 12).Value = 1890.6666
$ws.Cells.Item(107, 13).Value = 810.8667
$ws.Cells.Item(107, 14).Value = -5730.6666

$ws = $wb.Worksheets.Item("CRP")
$ws.Cells.Item(31, 8).Value = 2117.694
$ws.Cells.Item(31, 9).Value = 1785.7333
$ws.Cells.Item(31, 10).Value = 2491.15
$ws.Cells.Item(31, 11).Value = 1785.7333
$ws.Cells.Item(31, 12).Value = 2491.15
$ws.Cells.Item(31, 13).Value = -1490.7333
$ws.Cells.Item(31, 14).Value = -3081.15
$ws.Cells.Item(34, 8).Value = 2117.694
$ws.Cells.Item(34, 9).Value = 1785.7333
$ws.Cells.Item(34, 10).Value = 2491.15
$ws.Cells.Item(34, 11).Value = 1785.7333
$ws.Cells.Item(34, 12).Value = 2491.15
$ws.Cells.Item(34, 13).Value = -1583.7333
$ws.Cells.Item(34, 14).Value = -2895.15
$ws.Cells.Item(58, 8).Value = 1681.7441
$ws.Cells.Item(58, 9).Value = 1000.0345
$ws.Cells.Item(58, 10).Value = 3093.8572
$ws.Cells.Item(58, 11).Value = 1000.0345
$ws.Cells.Item(58, 12).Value = 3093.8572
$ws.Cells.Item(58, 13).Value = -797.0345
$ws.Cells.Item(58, 14).Value = -3499.8572
$ws.Cells.Item(132, 8).Value = 2047.8286
$ws.Cells.Item(132, 9).Value = 1547.4
$ws.Cells.Item(132, 10).Value = 5050.4
$ws.Cells.Item(132, 11).Value = 4642.200000000001
$ws.Cells.Item(132, 12).Value = 15151.2
$ws.Cells.Item(132, 13).Value = -2112.200000000001
$ws.Cells.Item(132, 14).Value = -20211.2
$ws.Cells.Item(134, 8).Value = 663176.1
$ws.Cells.Item(134, 9).Value = 1963.6562
$ws.Cells.Item(134, 10).Value = 5952876
$ws.Cells.Item(134, 11).Value = 5890.9686
$ws.Cells.Item(134, 12).Value = 17858628
$ws.Cells.Item(134, 13).Value = -3355.9686
$ws.Cells.Item(134, 14).Value = -17863698
$ws.Cells.Item(136, 8).Value = 1681.7441
$ws.Cells.Item(136, 9).Value = 1000.0345
$ws.Cells.Item(136, 10).Value = 3093.8572
$ws.Cells.Item(136, 11).Value = 3000.1035
$ws.Cells.Item(136, 12).Value = 9281.571599999999
$ws.Cells.Item(136, 13).Value = -450.1035000000002
$ws.Cells.Item(136, 14).Value = -14381.5716

$ws = $wb.Worksheets.Item("CUL")
$ws.Cells.Item(101, 8).Value = 21333.334
$ws.Cells.Item(101, 10).Value = 21333.334
$ws.Cells.Item(101, 12).Value = 64000.00199999999
$ws.Cells.Item(101, 14).Value = -68868.00199999999

$ws = $wb.Worksheets.Item("GSM")
$ws.Cells.Item(70, 8).Value = 36833.332
$ws.Cells.Item(70, 9).Value = 53250
$ws.Cells.Item(70, 10).Value = 4000
$ws.Cells.Item(70, 11).Value = 53250
$ws.Cells.Item(70, 12).Value = 4000
$ws.Cells.Item(70, 13).Value = -52980
$ws.Cells.Item(70, 14).Value = -4540
$ws.Cells.Item(73, 8).Value = 36833.332
$ws.Cells.Item(73, 9).Value = 53250
$ws.Cells.Item(73, 10).Value = 4000
$ws.Cells.Item(73, 11).Value = 53250
$ws.Cells.Item(73, 12).Value = 4000
$ws.Cells.Item(73, 13).Value = -52314
$ws.Cells.Item(73, 14).Value = -5872
$ws.Cells.Item(102, 8).Value = 3399.3076
$ws.Cells.Item(102, 9).Value = 4934.2856
$ws.Cells.Item(102, 10).Value = 1608.5
$ws.Cells.Item(102, 11).Value = 4934.2856
$ws.Cells.Item(102, 12).Value = 1608.5
$ws.Cells.Item(102, 13).Value = -3312.2856
$ws.Cells.Item(102, 14).Value = -4852.5
$ws.Cells.Item(122, 8).Value = 2684.25
$ws.Cells.Item(122, 9).Value = 2558.9092
$ws.Cells.Item(122, 11).Value = 7676.7276
$ws.Cells.Item(122, 13).Value = -5226.7276
$ws.Cells.Item(126, 8).Value = 3722.111
$ws.Cells.Item(126, 9).Value = 2399.7778
$ws.Cells.Item(126, 10).Value = 5044.4443
$ws.Cells.Item(126, 11).Value = 7199.3334
$ws.Cells.Item(126, 12).Value = 15133.3329
$ws.Cells.Item(126, 13).Value = -4729.3334
$ws.Cells.Item(126, 14).Value = -20073.3329

$ws = $wb.Worksheets.Item("LTW")
$ws.Cells.Item(136, 8).Value = 1559.5264
$ws.Cells.Item(136, 9).Value = 1575.2452
$ws.Cells.Item(136, 10).Value = 1351.25
$ws.Cells.Item(136, 11).Value = 4725.7356
$ws.Cells.Item(136, 12).Value = 4053.75
$ws.Cells.Item(136, 13).Value = -2175.7356
$ws.Cells.Item(136, 14).Value = -9153.75

$ws = $wb.Worksheets.Item("WVR")
$ws.Cells.Item(12, 8).Value = 1000
$ws.Cells.Item(12, 9).Value = 0
$ws.Cells.Item(12, 10).Value = 1000
$ws.Cells.Item(12, 11).Value = 0
$ws.Cells.Item(12, 12).Value = 1000
$ws.Cells.Item(12, 13).ClearContents()
$ws.Cells.Item(12, 14).Value = -1284
$ws.Cells.Item(136, 8).Value = 1396.7059
$ws.Cells.Item(136, 9).Value = 835.5897
$ws.Cells.Item(136, 10).Value = 3220.3333
$ws.Cells.Item(136, 11).Value = 2506.7691
$ws.Cells.Item(136, 12).Value = 9660.999899999999
$ws.Cells.Item(136, 13).Value = 43.23090000000002
$ws.Cells.Item(136, 14).Value = -14760.9999
